$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the "Absent" (H) column so every row has a numeric value:
# H4 was 0 -> should be 1 (absent, since D4 = 0 i.e. no attendance that day)
$ws.Range("H4").Value = 1

# H5 was blank (inline string) -> should be numeric 0 (present, D5 = 1)
$ws.Range("H5").Value = 0

# H9 was 0 -> should be 1 (absent, since D9 = 0)
$ws.Range("H9").Value = 1

# H10 was blank (inline string) -> should be numeric 0 (present, D10 = 1)
$ws.Range("H10").Value = 0
